# Fruta / hortaliza, semanal
# A new weekly observation is inserted at row 9 (Guayaba, Vega Modelo de
# Temuco), pushing the existing rows 9-24 down to rows 10-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 9, shifting rows 9:24 -> 10:25.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Cells.Item(9, 1).Value  = 10
$ws.Cells.Item(9, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(9, 3).Value  = "La Araucanía"
$ws.Cells.Item(9, 4).Value  = 45044
$ws.Cells.Item(9, 5).Value  = 9
$ws.Cells.Item(9, 6).Value  = "Fruta"
$ws.Cells.Item(9, 7).Value  = 100108
$ws.Cells.Item(9, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(9, 9).Value  = 100108001
$ws.Cells.Item(9, 10).Value = "Guayaba"
$ws.Cells.Item(9, 11).Value = "Sin especificar"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 150
$ws.Cells.Item(9, 14).Value = 3500
$ws.Cells.Item(9, 15).Value = 3500
$ws.Cells.Item(9, 16).Value = 3500
$ws.Cells.Item(9, 17).Value = "$/kilo"
$ws.Cells.Item(9, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(9, 19).Value = 3500
$ws.Cells.Item(9, 20).Value = 1
